$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- WebForm User Assignment execution ---
# Refresh the per-row phone-number ("RRN") values in column F, including a
# previously-missing value for row 5. These must stay TEXT (as the source
# data always has been) rather than be auto-coerced to numbers, so format
# the range as Text before writing, then drop back to the default style.

$phoneRange = $ws.Range("F2:F10")
$phoneRange.NumberFormat = "@"

$ws.Range("F2").Value = "9840059203"
$ws.Range("F3").Value = "9840016497"
$ws.Range("F4").Value = "9840036331"
$ws.Range("F5").Value = "9840047372"
$ws.Range("F6").Value = "9840027047"
$ws.Range("F7").Value = "9840003494"
$ws.Range("F8").Value = "9840079207"
$ws.Range("F9").Value = "9840060057"
$ws.Range("F10").Value = "9840063038"

$phoneRange.Style = "Normal"

# --- Swap the matched-user-position flags on row 2 (AM2/AN2) ---
$matchRange = $ws.Range("AM2:AN2")
$matchRange.NumberFormat = "@"

$ws.Range("AM2").Value = "0"
$ws.Range("AN2").Value = "2"

$matchRange.Style = "Normal"

# --- Update the active selection left on the sheet after the run ---
$ws.Range("AP12").Select()
